# Edit script: add 4 new Data Quality Assertion rows (duplicate_value_combination and
# geo_spatial_accuracy_precision) to the "Data quality assertion" and "Definition of
# assertions" sheets, inserted before the existing "date_format_validation" rows.

$wb = $excel.ActiveWorkbook

# --- Sheet: "Data quality assertion" ---
$ws2 = $wb.Worksheets.Item("Data quality assertion")
$ws2.Rows.Item(28).Resize(4).Insert()

$ws2.Range("A28").Value = 'duplicate_value_combination:inferred_duplicate'
$ws2.Range("B28").Value = ''
$ws2.Range("A29").Value = 'duplicate_value_combination:inferred_nonduplicate'
$ws2.Range("B29").Value = ''
$ws2.Range("A30").Value = 'geo_spatial_accuracy_precision:low_precision'
$ws2.Range("B30").Value = ''
$ws2.Range("A31").Value = 'geo_spatial_accuracy_precision:high_precision'
$ws2.Range("B31").Value = ''

# --- Sheet: "Definition of assertions" ---
$ws3 = $wb.Worksheets.Item("Definition of assertions")
$ws3.Rows.Item(28).Resize(4).Insert()

$ws3.Range("A28").Value = 'duplicate_value_combination:inferred_duplicate'
$ws3.Range("B28").Value = 'data_quality'
$ws3.Range("C28").Value = 'various_fields'
$ws3.Range("D28").Value = 'inferred_duplicate'
$ws3.Range("E28").Value = 'Indicates that the record has a combination of values across multiple fields that are identical to other records.'
$ws3.Range("F28").Value = 'If the record has a combination of values across the specified fields that is identical to other records in the dataset, label it as ''duplicate_combination''. This implies redundancy in data values for multiple records.'

$ws3.Range("A29").Value = 'duplicate_value_combination:inferred_nonduplicate'
$ws3.Range("B29").Value = 'data_quality'
$ws3.Range("C29").Value = 'various_fields'
$ws3.Range("D29").Value = 'inferred_nonduplicate'
$ws3.Range("E29").Value = 'Indicates that the record has a unique combination of values across multiple fields that is not shared by other records.'
$ws3.Range("F29").Value = 'If the record has a unique combination of values across the specified fields, label it as ''unique_combination''. This means that no other records share this exact combination.'

$ws3.Range("A30").Value = 'geo_spatial_accuracy_precision:low_precision'
$ws3.Range("B30").Value = 'geo'
$ws3.Range("C30").Value = 'geo:hasMetricSpatialAccuracy'
$ws3.Range("D30").Value = 'low_precision'
$ws3.Range("E30").Value = 'Indicates that the spatial accuracy is low, either because the value of coordinateUncertaintyInMeters is empty or exceeds 10,000 meters.'
$ws3.Range("F30").Value = 'If the ''coordinateUncertaintyInMeters'' field is empty or its value exceeds 10,000, label the record as ''low_precision''. This indicates that the precision of the spatial accuracy is insufficient.'

$ws3.Range("A31").Value = 'geo_spatial_accuracy_precision:high_precision'
$ws3.Range("B31").Value = 'geo'
$ws3.Range("C31").Value = 'geo:hasMetricSpatialAccuracy'
$ws3.Range("D31").Value = 'high_precision'
$ws3.Range("E31").Value = 'Indicates that the spatial accuracy is high, meaning the value of coordinateUncertaintyInMeters is less than or equal to 10,000 meters.'
$ws3.Range("F31").Value = 'If the ''coordinateUncertaintyInMeters'' field contains a value of 10,000 or less, label the record as ''high_precision''. This indicates that the precision of the spatial accuracy is adequate.'

Write-Host "Done inserting new assertion rows."
